# "add discpritions in sprint backlog"
# Adds task descriptions (column D) for three Sprint Backlog rows, fixes a
# typo in an existing description, fills in the "Effort Plan Updated (h)"
# column (J) to mirror the "Effort Plan Original (h)" column (I) for the
# rows that previously had no planned/updated effort recorded, and adjusts
# the row heights that Excel recomputed as a result of the new wrapped text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sprint Backlog")

# --- New task descriptions (column D) ---------------------------------
$ws.Range("D10").Value = "Create the Java patient model out of the DB "
$ws.Range("D15").Value = "Creating a model for the Patient Info (may contain more informatino than just the patient model)"
$ws.Range("D17").Value = "Create the presenter for the Patient Info (acts between View Interface and Model)"

# --- Typo fix in an existing description -------------------------------
$ws.Range("D19").Value = "Shows the objective (they can be created, opened or changed)"

# --- Fill "Effort Plan Updated (h)" (column J) to match column I --------
$ws.Range("J9").Value = 5
$ws.Range("J10").Value = 3
$ws.Range("J13").Value = 9
$ws.Range("J15").Value = 4
$ws.Range("J16").Value = 4
$ws.Range("J17").Value = 3
$ws.Range("J18").Value = 5
$ws.Range("J19").Value = 9
$ws.Range("J20").Value = 4
$ws.Range("J21").Value = 6

# --- Row heights recomputed by Excel after the wrapped text changed ----
$ws.Rows.Item(9).RowHeight = 60
$ws.Rows.Item(10).RowHeight = 30
$ws.Rows.Item(15).RowHeight = 75
$ws.Rows.Item(17).RowHeight = 60
$ws.Rows.Item(19).RowHeight = 45
$ws.Rows.Item(20).RowHeight = 30
